$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of price reports was collected (commit: "Fruta / hortaliza, semanal").
# The new observations are prepended above the existing data block, pushing the
# previously-existing rows 109-132 down to 111-134.
$ws.Rows("109:110").Insert()

# Row 109: Feria Lagunitas de Puerto Montt - Mango, "Primera" quality, week of 2021-12-21
$ws.Range("A109").Value = 4
$ws.Range("B109").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C109").Value = "Los Lagos"
$ws.Range("D109").Value = "2021-12-21"
$ws.Range("E109").Value = 10
$ws.Range("F109").Value = "Fruta"
$ws.Range("G109").Value = 100108
$ws.Range("H109").Value = "Tropicales y subtropicales"
$ws.Range("I109").Value = 100108002
$ws.Range("J109").Value = "Mango"
$ws.Range("K109").Value = "Sin especificar"
$ws.Range("L109").Value = "Primera"
$ws.Range("M109").Value = 200
$ws.Range("N109").Value = 8500
$ws.Range("O109").Value = 9000
$ws.Range("P109").Value = 8750
$ws.Range("Q109").Value = "$/bandeja 4 kilos"
$ws.Range("R109").Value = "Perú"
$ws.Range("S109").Value = 2188
$ws.Range("T109").Value = 4

# Row 110: Feria Lagunitas de Puerto Montt - Mango, "Segunda" quality, week of 2021-12-21
$ws.Range("A110").Value = 4
$ws.Range("B110").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C110").Value = "Los Lagos"
$ws.Range("D110").Value = "2021-12-21"
$ws.Range("E110").Value = 10
$ws.Range("F110").Value = "Fruta"
$ws.Range("G110").Value = 100108
$ws.Range("H110").Value = "Tropicales y subtropicales"
$ws.Range("I110").Value = 100108002
$ws.Range("J110").Value = "Mango"
$ws.Range("K110").Value = "Sin especificar"
$ws.Range("L110").Value = "Segunda"
$ws.Range("M110").Value = 100
$ws.Range("N110").Value = 6000
$ws.Range("O110").Value = 6000
$ws.Range("P110").Value = 6000
$ws.Range("Q110").Value = "$/bandeja 4 kilos"
$ws.Range("R110").Value = "Perú"
$ws.Range("S110").Value = 1500
$ws.Range("T110").Value = 4
